$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 337, shifting the existing rows 337:392 down to 338:393
$ws.Rows.Item(337).Insert()

# Populate the newly inserted row 337 with the new weekly price record
$ws.Range("A337").Value = 8
$ws.Range("B337").Value = "Terminal La Palmera de La Serena"
$ws.Range("C337").Value = "Coquimbo"
$ws.Range("D337").Value = (Get-Date -Year 2023 -Month 5 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E337").Value = 4
$ws.Range("F337").Value = 100112031
$ws.Range("G337").Value = "Poroto verde"
$ws.Range("H337").Value = "Magnum"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 200
$ws.Range("K337").Value = 22000
$ws.Range("L337").Value = 23000
$ws.Range("M337").Value = 22500
$ws.Range("N337").Value = '$/malla 25 kilos'
$ws.Range("O337").Value = "Provincia de Limarí"
$ws.Range("P337").Value = 900
$ws.Range("Q337").Value = 25
$ws.Range("R337").Value = "Hortaliza"
